# Fix for "interação com input não clicável": the first cell (A1) was an
# empty input-looking cell that wasn't clickable/interactive; give it the
# text "teste" so it actually holds content, and leave the selection on
# the cell the user ends up on (E3) after interacting with the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 gets the literal text "teste" (was empty before).
$ws.Range("A1").Value = "teste"

# Move/leave the active selection on E3, matching the saved workbook state.
$ws.Range("E3").Select()
